# Atividade de controle de versoes.pptx
# Commit: "Corrigido posicionamento de dados do alunos Antonio"
#   (Fixed positioning of student Antonio's data)
#
# Slide 4 ("Aluno" / Antonio Carlos Gomes Tabosa) had his photo placeholder
# sitting on top of / too far right over his text, and his name line was a
# single run "Nome:Antônio Carlos Gomes Tabosa" (missing the space after
# the colon). This reflows the name into its separate pieces and moves the
# photo placeholder to its corrected position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Fix the "Nome:" line -------------------------------------------------
# TextBox 5 is the shape holding Nome / Email / Foto.
$txBox = $s.Shapes.Item(1)

$nameLine = $txBox.TextFrame.TextRange.Paragraphs(1, 1)

# Replace the run with "Nome" then grow it with the remaining pieces so the
# line ends up as three runs: "Nome" + ": Antônio " + "Carlos Gomes Tabosa"
# (note the fixed spacing: "Nome: Antônio ..." instead of "Nome:Antônio ...").
$nameLine.Text = "Nome"
$rest = $nameLine.InsertAfter(": Antônio ")
$rest.InsertAfter("Carlos Gomes Tabosa") | Out-Null

# Editing the text re-flows the auto-fit textbox; the box itself did not
# actually move or resize in the original edit, so put its height back to
# its original 4030980 EMU (317.4pt).
$txBox.Height = 317.40001

# --- Reposition Antonio's photo placeholder -------------------------------
$photo = $s.Shapes.Item(3)
$photo.Left = 3305166 / 914400 * 72
$photo.Top = 3506170 / 914400 * 72
